# The commit reorders the weekly price-report rows (rows 2-17) of the
# single data sheet: each row's full set of values (date, volume,
# min/max/avg price, unit, origin, etc.) moves to a different row
# position, while rows 12 and 16 stay put. This reproduces that
# permutation by first snapshotting every source row (as a whole, via a
# single Range.Value read) and then writing each snapshot back into its
# new destination row, so that overlapping reads/writes never clobber
# data that still needs to be used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current contents of rows 2..17 (columns A..T) before
# writing anything, since the reassignment below is a permutation and
# some rows depend on others' original values.
$rows = @{}
for ($r = 2; $r -le 17; $r++) {
    $rows[$r] = $ws.Range("A$r`:T$r").Value()
}

# Maps destination row -> source row (i.e. destination row ends up with
# the values that source row used to hold).
$mapping = @{
    2  = 17
    3  = 5
    4  = 11
    5  = 6
    6  = 13
    7  = 15
    8  = 2
    9  = 3
    10 = 4
    11 = 14
    12 = 12
    13 = 8
    14 = 7
    15 = 10
    16 = 16
    17 = 9
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $ws.Range("A$destRow`:T$destRow").Value = $rows[$srcRow]
}
